# Insert the new "2022-Q3" worksheet (a copy of the existing "2022-Q2"
# worksheet, positioned right before it) and populate it with the new
# quarter's fund-holding data, then update the "总计" (totals) sheet with a
# new row summarising the 2022-Q3 quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the "2022-Q3" sheet by duplicating "2022-Q2" (this keeps all
#    sheet-level formatting / sheetPr / column styles identical to its
#    siblings) and inserting the copy immediately before "2022-Q2".
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2, $null)

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# ---------------------------------------------------------------------
# 2. Overwrite the data in "2022-Q3" with the new quarter's fund table.
#    Columns B-G mirror the source data's own convention of storing
#    these figures as text, so we briefly force a text number format
#    while assigning them, then drop back to General (keeps the cells
#    free of any lingering explicit number format).
# ---------------------------------------------------------------------
$q3.Range("B2:G9").NumberFormat = "@"

$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "000118"
$q3.Range("C2").Value = "广发聚鑫债券A"
$q3.Range("D2").Value = "145.84"
$q3.Range("E2").Value = "20.01"
$q3.Range("F2").Value = "1.00"
$q3.Range("G2").Value = "1.4584"
$q3.Range("H2").Value = 9

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "009121"
$q3.Range("C3").Value = "广发招享混合A"
$q3.Range("D3").Value = "58.31"
$q3.Range("E3").Value = "23.12"
$q3.Range("F3").Value = "1.08"
$q3.Range("G3").Value = "0.6297"
$q3.Range("H3").Value = 8

$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "260103"
$q3.Range("C4").Value = "景顺长城动力平衡混合"
$q3.Range("D4").Value = "10.57"
$q3.Range("E4").Value = "68.77"
$q3.Range("F4").Value = "2.66"
$q3.Range("G4").Value = "0.2812"
$q3.Range("H4").Value = 9

$q3.Range("A5").Value = 3
$q3.Range("B5").Value = "013880"
$q3.Range("C5").Value = "广发招享混合C"
$q3.Range("D5").Value = "25.35"
$q3.Range("E5").Value = "23.12"
$q3.Range("F5").Value = "1.08"
$q3.Range("G5").Value = "0.2738"
$q3.Range("H5").Value = 8

$q3.Range("A6").Value = 4
$q3.Range("B6").Value = "000119"
$q3.Range("C6").Value = "广发聚鑫债券C"
$q3.Range("D6").Value = "16.95"
$q3.Range("E6").Value = "20.01"
$q3.Range("F6").Value = "1.00"
$q3.Range("G6").Value = "0.1695"
$q3.Range("H6").Value = 9

$q3.Range("A7").Value = 5
$q3.Range("B7").Value = "010949"
$q3.Range("C7").Value = "景顺长城研究驱动三年持有期混合"
$q3.Range("D7").Value = "2.48"
$q3.Range("E7").Value = "67.69"
$q3.Range("F7").Value = "2.67"
$q3.Range("G7").Value = "0.0662"
$q3.Range("H7").Value = 9

$q3.Range("A8").Value = 6
$q3.Range("B8").Value = "011404"
$q3.Range("C8").Value = "融通鑫新成长混合C"
$q3.Range("D8").Value = "1.75"
$q3.Range("E8").Value = "94.07"
$q3.Range("F8").Value = "2.93"
$q3.Range("G8").Value = "0.0513"
$q3.Range("H8").Value = 9

$q3.Range("A9").Value = 7
$q3.Range("B9").Value = "011403"
$q3.Range("C9").Value = "融通鑫新成长混合A"
$q3.Range("D9").Value = "0.39"
$q3.Range("E9").Value = "94.07"
$q3.Range("F9").Value = "2.93"
$q3.Range("G9").Value = "0.0114"
$q3.Range("H9").Value = 9

$q3.Range("B2:G9").ClearFormats()

# Restore the column-A / header styling that ClearFormats just wiped by
# re-copying it across from the sibling sheet (keeps the same style
# index instead of minting a new one).
$q2.Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$q2.Range("A2:A9").Copy()
$q3.Range("A2:A9").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3. Update the "总计" (totals) sheet: shift the existing quarters down
#    one row and insert the new 2022-Q3 totals as row 2.
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

$totals.Range("A8").Value = 6
$totals.Range("B8").Value = "2021-Q1"
$totals.Range("C8").Value = 2
$totals.Range("D8").Value = 0.39

$totals.Range("A7").Value = 5
$totals.Range("B7").Value = "2021-Q2"
$totals.Range("C7").Value = 5
$totals.Range("D7").Value = 1.56

$totals.Range("A6").Value = 4
$totals.Range("B6").Value = "2021-Q3"
$totals.Range("C6").Value = 7
$totals.Range("D6").Value = 1.58

$totals.Range("A5").Value = 3
$totals.Range("B5").Value = "2021-Q4"
$totals.Range("C5").Value = 10
$totals.Range("D5").Value = 2.27

$totals.Range("A4").Value = 2
$totals.Range("B4").Value = "2022-Q1"
$totals.Range("C4").Value = 24
$totals.Range("D4").Value = 5.27

$totals.Range("A3").Value = 1
$totals.Range("B3").Value = "2022-Q2"
$totals.Range("C3").Value = 8
$totals.Range("D3").Value = 1.01

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q3"
$totals.Range("C2").Value = 8
$totals.Range("D2").Value = 2.94

# A8 is a brand new cell - copy the column's style (index "2") onto it
# instead of leaving it in the default style.
$totals.Range("A3").Copy()
$totals.Range("A8").PasteSpecial(-4122)
